$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price / "Volume(1h)" values for this scrape run
$ws.Range("D2").Value = "25.847.47"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.740.88"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'224.55"
$ws.Range("E5").Value = "  -5.39%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5150"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.2799"
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("D9").Value = "'38.99"
$ws.Range("E9").Value = "  -4.87%  "
$ws.Range("D10").Value = "'0.06079"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "1.739.64"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'0.06952"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "'0.6314"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").Value = "'4.483"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "'76.25"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "25.862.36"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'11.39"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "'0.000006550"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "'4.069"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'8.403"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("D25").Value = "'5.097"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'137.74"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").Value = "'1.814"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'14.93"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "'102.29"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'0.08258"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'3.605"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "'3.396"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'0.04380"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'2.623"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "'0.9631"
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("D37").Value = "'0.5988"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'2.667"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'1.895"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "'100.61"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").Value = "'0.3809"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.7215"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").Value = "'4.884"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'0.05454"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'6.255"
$ws.Range("E47").Value = "  +5.31%  "
$ws.Range("D48").Value = "'0.1092"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'52.06"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'29.61"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  -0.07%  "

# Plain numeric-looking new values (e.g. "224.55", "0.000006550") were
# entered above with a leading apostrophe so Excel keeps them as exact
# text (preserving trailing/leading zeros and the fixed layout the sheet
# relies on) instead of silently re-parsing them as numbers. That leaves
# a "quote prefix" flag on each affected cell's style, so clear it back
# to the default worksheet style, cell by cell (a multi-area Range.Style
# assignment only touches the first area).
$textCells = @(
    "D5", "D7", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
